# Prefab System Tutorial Renewal
# The seventh tutorial slide ("일곱 번째 튜토리얼" / Prefab Release Button walkthrough)
# and its embedded resource images are removed from the deck.

$p = $ppt.ActivePresentation

# The seventh slide (sldId 294, slides/slide7.xml) is the last slide in the deck.
# Deleting it also removes its relationship entry, embedded picture parts that
# are no longer referenced, and its <p:sldId> entry from the presentation's
# slide id list.
$s = $p.Slides.Item($p.Slides.Count)
$s.Delete()
